$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45171 = 2023-09-02) for
# every data row (2..236). The commit bumps that date by one day to 45172
# (2023-09-03) across the whole column.
$ws.Range("C2:C236").Value = 45172
